# Katpally_LabExam03Grading.xlsx - grading pass
# "kalyankar to pusapati done"
#
# Fill in the "Points for grading" (column E) scores for the
# CustomerMapping Class and CustomerMappingDriver Class sections,
# matching the max points already recorded in column D for those rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# CustomerMapping Class block (rows 3-6)
$ws.Range("E3").Value = 1
$ws.Range("E4").Value = 2
$ws.Range("E5").Value = 2
$ws.Range("E6").Value = 2

# CustomerMappingDriver Class block (rows 10-14)
$ws.Range("E10").Value = 2
$ws.Range("E11").Value = 2
$ws.Range("E12").Value = 2
$ws.Range("E13").Value = 2
$ws.Range("E14").Value = 2

# Leave the cursor/selection where grading left off.
$ws.Range("E15").Select()
